$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.914.22"
$ws.Range("E2").Value = "  -4.09%  "
$ws.Range("D3").Value = "1.741.60"
$ws.Range("E3").Value = "  -4.50%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "225.05"
$ws.Range("E5").Value = "  -3.92%  "
$ws.Range("D6").Value = "0.5764"
$ws.Range("E6").Value = "  -3.76%  "
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.2695"
$ws.Range("E8").Value = "  -1.96%  "
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").Value = "0.06555"
$ws.Range("E10").Value = "  -5.49%  "
$ws.Range("D11").Value = "0.07497"
$ws.Range("E11").Value = "  -1.56%  "
$ws.Range("D12").Value = "1.739.90"
$ws.Range("E12").Value = "  -4.56%  "
$ws.Range("D13").Value = "4.680"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").Value = "0.6006"
$ws.Range("E14").Value = "  -3.73%  "
$ws.Range("D15").Value = "1.977.70"
$ws.Range("E15").Value = "  -4.50%  "
$ws.Range("D16").Value = "73.48"
$ws.Range("E16").Value = "  -5.19%  "
$ws.Range("D17").Value = "0.000008579"
$ws.Range("E17").Value = "  -11.07%  "
$ws.Range("D18").Value = "27.938.15"
$ws.Range("E18").Value = "  -2.79%  "
$ws.Range("D19").Value = "5.280"
$ws.Range("E19").Value = "  -5.16%  "
$ws.Range("D20").Value = "1.005"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").Value = "203.47"
$ws.Range("E21").Value = "  -6.03%  "
$ws.Range("D22").Value = "11.23"
$ws.Range("E22").Value = "  -2.38%  "
$ws.Range("D23").Value = "6.609"
$ws.Range("E23").Value = "  -3.48%  "
$ws.Range("D24").Value = "1.006"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "150.22"
$ws.Range("E25").Value = "  -4.03%  "
$ws.Range("D26").Value = "7.960"
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("D27").Value = "0.1226"
$ws.Range("E27").Value = "  -4.24%  "
$ws.Range("D28").Value = "15.99"
$ws.Range("E28").Value = "  -2.74%  "
$ws.Range("D29").Value = "1.390"
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("D30").Value = "0.06017"
$ws.Range("E30").Value = "  -5.97%  "
$ws.Range("D31").Value = "1.380"
$ws.Range("E31").Value = "  -3.98%  "
$ws.Range("D32").Value = "3.708"
$ws.Range("E32").Value = "  -3.06%  "
$ws.Range("D33").Value = "3.699"
$ws.Range("E33").Value = "  -1.34%  "
$ws.Range("D34").Value = "1.669"
$ws.Range("E34").Value = "  -3.19%  "
$ws.Range("D35").Value = "1.029"
$ws.Range("E35").Value = "  -5.21%  "
$ws.Range("D36").Value = "0.6293"
$ws.Range("E36").Value = "  -2.34%  "
$ws.Range("D37").Value = "2.434"
$ws.Range("E37").Value = "  -4.10%  "
$ws.Range("D38").Value = "2.645"
$ws.Range("E38").Value = "  -3.64%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01657"
$ws.Range("E39").Value = "  -5.28%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "6.228"
$ws.Range("E40").Value = "  -4.92%  "
$ws.Range("D41").Value = "1.120.09"
$ws.Range("E41").Value = "  -2.43%  "
$ws.Range("D42").Value = "0.8602"
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("D43").Value = "1.008"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("D44").Value = "99.02"
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("D45").Value = "1.889.00"
$ws.Range("E45").Value = "  -4.32%  "
$ws.Range("D46").Value = "58.74"
$ws.Range("E46").Value = "  -4.87%  "
$ws.Range("D47").Value = "0.00000000108"
$ws.Range("E47").Value = "  -4.64%  "
$ws.Range("E48").Value = "  -2.83%  "
$ws.Range("D49").Value = "8.216"
$ws.Range("E49").Value = "  -2.38%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.05389"
$ws.Range("E50").Value = "  -2.14%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "0.4447"
$ws.Range("E51").Value = "  -2.01%  "
